# tabla_gaussSeidel: refresh the iteration table with the Gauss-Seidel
# root-approximation results (grows from 2 to 15 iterations / rows 2-16)
# and fix the xn / error values for the first two iterations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values in this table (including the numeric-looking ones such as
# "1.0", "0.25", "3") are stored as plain TEXT, matching how the sheet was
# originally authored (every data cell is an inline/shared string, never a
# real number). Writing through `.Value` would let Excel auto-convert
# numeric-looking strings into actual numbers, so instead we write each
# value as a `="..."` text formula and then "Paste Values" it back onto
# itself; that bakes in the literal text without touching NumberFormat (and
# therefore without minting any new cell style).
function Set-TextValue($cell, [string]$value) {
    $escaped = $value.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
}

# Row 2 (iteration 1) -- B2/C2 change
Set-TextValue $ws.Cells.Item(2, 2) "[0;0.5]"
Set-TextValue $ws.Cells.Item(2, 3) "inf"

# Row 3 (iteration 2) -- B3/C3 change
Set-TextValue $ws.Cells.Item(3, 2) "[0.333333333333333;0.333333333333333]"
Set-TextValue $ws.Cells.Item(3, 3) "1.0"

# Rows 4-16 (iterations 3-15) -- newly added
$rows = @(
    @(4,  "3",  "[0.444444444444444;0.277777777777778]", "0.25"),
    @(5,  "4",  "[0.481481481481481;0.259259259259259]", "0.0769230769230769"),
    @(6,  "5",  "[0.493827160493827;0.253086419753086]", "0.025"),
    @(7,  "6",  "[0.497942386831276;0.251028806584362]", "0.0082644628099173"),
    @(8,  "7",  "[0.499314128943759;0.250342935528121]", "0.0027472527472527"),
    @(9,  "8",  "[0.499771376314586;0.250114311842707]", "0.0009149130832571"),
    @(10, "9",  "[0.499923792104862;0.250038103947569]", "0.0003048780487804"),
    @(11, "10", "[0.499974597368287;0.250012701315856]", "0.0001016156894624"),
    @(12, "11", "[0.499991532456096;0.250004233771952]", "3.3870749220958e-05"),
    @(13, "12", "[0.499997177485365;0.250001411257317]", "1.12901222720193e-05"),
    @(14, "13", "[0.499999059161788;0.250000470419106]", "3.76335992777886e-06"),
    @(15, "14", "[0.499999686387263;0.250000156806369]", "1.25445173560849e-06"),
    @(16, "15", "[0.499999895462421;0.25000005226879]", "4.18150403686329e-07")
)

foreach ($entry in $rows) {
    $r = $entry[0]
    Set-TextValue $ws.Cells.Item($r, 1) $entry[1]
    Set-TextValue $ws.Cells.Item($r, 2) $entry[2]
    Set-TextValue $ws.Cells.Item($r, 3) $entry[3]
}

# Bake the `="..."` formulas above into plain literal text values.
$dataRange = $ws.Range("A2:C16")
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
